$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (pushes current rows 5,6 down to 6,7)
$ws.Rows.Item(5).Insert()

# Update existing rows with new values (B column unchanged where given)
# Row 2: ln_pgdp
$ws.Range("C2").Value = 0.3365521852103142
$ws.Range("D2").Value = 99.1231071854279
$ws.Range("G2").Value = 0.09462444644427603
$ws.Range("H2").Value = 0.9246191308781114

# Row 3: ln_pop_density
$ws.Range("C3").Value = 3.271726598574773
$ws.Range("D3").Value = 90.90949586190925
$ws.Range("G3").Value = 0.9198731486877424
$ws.Range("H3").Value = 0.3577092649888407

# Row 4: tertiary_share
$ws.Range("C4").Value = 20.489872325343
$ws.Range("D4").Value = 54.67191585747635
$ws.Range("G4").Value = 5.76089804702315
$ws.Range("H4").Value = [double]"9.170212484141563e-09"

# Row 5 (new): tertiary_share_sq
$ws.Range("A5").Value = "tertiary_share_sq"
$ws.Range("B5").Value = 45.61125931829044
$ws.Range("C5").Value = 21.53136706643642
$ws.Range("D5").Value = 52.79374569295835
$ws.Range("E5").Value = 13.56053058748441
$ws.Range("F5").Value = [double]"7.429898583901447e-41"
$ws.Range("G5").Value = 6.053722957041203
$ws.Range("H5").Value = [double]"1.581736069367004e-09"
$ws.Range("I5").Value = "需检查"

# Row 6 (formerly row 5): ln_fdi - full row content changed
$ws.Range("A6").Value = "ln_fdi"
$ws.Range("B6").Value = 32.66295094300895
$ws.Range("C6").Value = -2.979673132347963
$ws.Range("D6").Value = 109.1224860164869
$ws.Range("E6").Value = 9.656641617854888
$ws.Range("F6").Value = [double]"8.645295030347428e-22"
$ws.Range("G6").Value = -0.8377598872433861
$ws.Range("H6").Value = 0.4022289290092775
$ws.Range("I6").Value = "OK"

# Row 7 (formerly row 6): ln_road_area - values changed
$ws.Range("A7").Value = "ln_road_area"
$ws.Range("B7").Value = -22.93534586918614
$ws.Range("C7").Value = -0.7723781488215303
$ws.Range("D7").Value = -96.63236755518378
$ws.Range("E7").Value = -6.731473306095445
$ws.Range("F7").Value = [double]"1.957702669230089e-11"
$ws.Range("G7").Value = -0.2171605414839901
$ws.Range("H7").Value = 0.8280972535816009
$ws.Range("I7").Value = "OK"
